$d = $word.ActiveDocument

# Locate the paragraph that ends the document's last list item
# ("And commit to it(save)") so the two new paragraphs land right after it.
$commitP = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*And commit to it(save)*") {
        $commitP = $d.Paragraphs.Item($i)
        break
    }
}

# The very last paragraph in the document is a bare, formatting-less
# paragraph (no pStyle / numbering). Inserting new paragraphs immediately
# before it (rather than cloning the preceding list paragraph) gives us a
# clean <w:p> with no inherited pStyle/numPr, matching what's needed here.
$lastP = $d.Paragraphs.Last
$lastP.Range.InsertParagraphBefore()
$lastP.Range.InsertParagraphBefore()

$newP1 = $d.Paragraphs.Item($commitP.Index + 1)
$newP2 = $d.Paragraphs.Item($commitP.Index + 2)

$newP1.Range.Text = "Task2"

foreach ($p in @($newP1, $newP2)) {
    $p.Shading.Texture = 0
    $p.Shading.ForegroundPatternColor = -16777216
    $p.Shading.BackgroundPatternColor = 16777215
    $p.SpaceBefore = 5
    $p.SpaceBeforeAuto = -1
    $p.SpaceAfter = 5
    $p.SpaceAfterAuto = -1
    $p.LineSpacingRule = 0
    $p.LineSpacing = 12
    $p.OutlineLevel = 3
}
